# Ahold.xlsx update: replace the serial-date "refilldate" values in the
# Product sheet (column K, rows 2-25) with free-text date strings.
# "20-Nov" (serial 43424) becomes "20 th November" and
# "30-Dec" (serial 43464) becomes "30 th December".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")

# Ordered row -> replacement text list (processed top to bottom so the
# shared-string table is populated in the same order as the source edit).
$replacements = @(
    @{ Row = 2;  Text = "20 th November" }
    @{ Row = 3;  Text = "30 th December" }
    @{ Row = 4;  Text = "30 th December" }
    @{ Row = 5;  Text = "30 th December" }
    @{ Row = 6;  Text = "20 th November" }
    @{ Row = 7;  Text = "30 th December" }
    @{ Row = 8;  Text = "30 th December" }
    @{ Row = 9;  Text = "30 th December" }
    @{ Row = 10; Text = "30 th December" }
    @{ Row = 11; Text = "30 th December" }
    @{ Row = 12; Text = "30 th December" }
    @{ Row = 13; Text = "30 th December" }
    @{ Row = 14; Text = "30 th December" }
    @{ Row = 15; Text = "30 th December" }
    @{ Row = 16; Text = "30 th December" }
    @{ Row = 17; Text = "30 th December" }
    @{ Row = 18; Text = "30 th December" }
    @{ Row = 19; Text = "30 th December" }
    @{ Row = 20; Text = "30 th December" }
    @{ Row = 21; Text = "30 th December" }
    @{ Row = 22; Text = "30 th December" }
    @{ Row = 23; Text = "30 th December" }
    @{ Row = 24; Text = "30 th December" }
    @{ Row = 25; Text = "30 th December" }
)

foreach ($entry in $replacements) {
    $ws.Cells.Item($entry.Row, 11).Value = $entry.Text
}

# The workbook was last touched with the selection sitting on K25.
$ws.Range("K25").Select()
